$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 13 ("testFolio", 1, YES).
# This pushes that row down to row 14 and creates a blank row 13
# that inherits formatting from the row above it (row 12).
$ws.Rows("13:13").Insert()

# New row 13: testScreening, 2, YES
$ws.Range("A13").Value = "testScreening"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "YES"

# Append a new row 15 with the same formatting as row 14 (testFolio, 1, YES),
# then set its values to testFolio, 2, YES.
$ws.Range("A14:C14").Copy()
[void]$ws.Range("A15:C15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A15").Value = "testFolio"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = "YES"

# Match the saved selection state from the workbook.
[void]$ws.Range("B13").Select()
